# Applies crypto price/volume updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.156.96"
$ws.Range("E2").Value = "  -2.32%  "
$ws.Range("D3").Value = "1.852.12"
$ws.Range("E3").Value = "  -1.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.27"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6834"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -5.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07667"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3035"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.10"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -6.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08135"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.76%  "
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7219"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.11%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.824.85"
$ws.Range("E13").Value = "  -3.36%  "
$ws.Range("E14").Value = "  -2.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.23"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.57%  "
$ws.Range("D16").Value = "29.139.24"
$ws.Range("E16").Value = "  -2.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007799"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.709"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -5.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.17"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "233.35"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -5.69%  "
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").Value = "2.098.54"
$ws.Range("E22").Value = "  -1.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.416"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.71"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.945"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.36%  "
$ws.Range("E27").Value = "  -5.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.02"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.953"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.393"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.508"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.484"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.003"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05155"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7021"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.023"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.03%  "
$ws.Range("E38").Value = "  -0.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01840"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.677"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9079"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.14%  "
$ws.Range("D42").Value = "1.101.42"
$ws.Range("E42").Value = "  +5.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.970"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4267"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -4.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "69.82"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.50%  "
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.46"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.767"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.53%  "
$ws.Range("D49").Value = "1.993.57"
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.127"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.918"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -7.41%  "
